# Updated formal attendance record
# - Renames table columns H ("Column1") and I ("Column2") to the new
#   meeting dates "2/2/2010" and "17/2/20102" respectively.
# - Fills in attendance for the new "2/2/2010" / "17/2/20102" meetings
#   for every attendee row (rows 3-8): "?" for the 2/2/2010 column and
#   "x" for the 17/2/20102 column (Haz gets "-" for 17/2/20102).
# - Leaves selection on J3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1) -----------------------------------------------
# Column I ("17/2/20102") is not a parseable date, so Excel keeps it as
# plain text automatically, but it still carries the same date display
# format as the other header cells.
$i1 = $ws.Cells.Item(1, 9)
$i1.Value = "17/2/20102"
$ws.Cells.Item(1, 2).Copy()
$i1.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column H ("2/2/2010") looks like a real date, so a plain .Value
# assignment on a General-formatted cell would get auto-converted to a
# date serial number. Force it to stay literal text first …
$h1 = $ws.Cells.Item(1, 8)
$h1.NumberFormat = "@"
$h1.Value = "2/2/2010"
# … then copy the existing date display format from the neighbouring
# header cell (B1) so H1 ends up visually formatted the same way as the
# other date headers, without disturbing the literal text value.
$ws.Cells.Item(1, 2).Copy()
$h1.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 3-8) ----------------------------------------------
$rows = 3..8
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = "?"
    $ws.Cells.Item($r, 9).Value = "x"
}
# Haz (row 6) didn't attend the 17/2/20102 meeting.
$ws.Cells.Item(6, 9).Value = "-"

# --- Selection ----------------------------------------------------------
$ws.Range("J3").Select()
